$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update status text from "In Translation" to "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Update timestamps
$overview.Range("G2").Value = "2016-08-21 11:02:46"
$dede.Range("H2").Value = "2016-08-21 11:02:46"
$zhcn.Range("H2").Value = "2016-08-21 11:02:42"

# Widen columns to fit the new, longer status text (matches Excel's computed
# autofit width of 17.2159881591797 for the OOXML "width" attribute, expressed
# here in COM "ColumnWidth" units: stored_width = ColumnWidth + 5/6)
$overview.Range("E:F").ColumnWidth = 16.333333333333332
$zhcn.Range("C:C").ColumnWidth = 16.333333333333332
$dede.Range("C:C").ColumnWidth = 16.333333333333332
